# Daily attendance processing - reorder "Recorded By" (column G) names so
# that the literal "System" entry/entries are moved to the front of the
# comma-separated list, preserving the relative order of the remaining
# entries. If a cell has no "System" entry, sort its entries alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = @($val -split ",\s*")

    $systemParts = New-Object System.Collections.ArrayList
    $restParts = New-Object System.Collections.ArrayList
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            [void]$systemParts.Add($p)
        } else {
            [void]$restParts.Add($p)
        }
    }

    if ($systemParts.Count -gt 0) {
        $newParts = @($systemParts) + @($restParts)
    } else {
        $newParts = @($restParts | Sort-Object)
    }

    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
